$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.542.73"
$ws.Range("E2").Value = "  +2.33%  "

$ws.Range("D3").Value = "1.911.13"
$ws.Range("E3").Value = "  +5.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.53"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5054"
$ws.Range("E7").Value = "  +1.66%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3953"
$ws.Range("E8").Value = "  +0.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09758"
$ws.Range("E9").Value = "  +1.97%  "

$ws.Range("E10").Value = "  +5.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.87"
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.551"
$ws.Range("E12").Value = "  +2.12%  "

$ws.Range("E13").Value = "  +3.85%  "

$ws.Range("D14").Value = "1.912.79"
$ws.Range("E14").Value = "  +5.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.564"

$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("E17").Value = "  +1.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.94"
$ws.Range("E18").Value = "  +2.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06665"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.11"
$ws.Range("E20").Value = "  +5.99%  "

$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.303"
$ws.Range("E22").Value = "  +6.80%  "

$ws.Range("D23").Value = "28.595.28"
$ws.Range("E23").Value = "  +2.31%  "

$ws.Range("E24").Value = "  +3.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.281"
$ws.Range("E25").Value = "  +1.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.760"
$ws.Range("E26").Value = "  +16.01%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.402"
$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.127.57"
$ws.Range("E28").Value = "  +5.33%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "21.27"
$ws.Range("E29").Value = "  +4.06%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "159.19"
$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.82"
$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.103"
$ws.Range("E32").Value = "  +7.04%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1072"
$ws.Range("E33").Value = "  +1.04%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.745"
$ws.Range("E34").Value = "  +3.71%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.644"
$ws.Range("E35").Value = "  +0.66%  "

$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.850"
$ws.Range("E36").Value = "  +10.96%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06780"
$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02444"
$ws.Range("E38").Value = "  +5.32%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.273"
$ws.Range("E39").Value = "  +9.93%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2233"
$ws.Range("E40").Value = "  +5.02%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.77"
$ws.Range("E41").Value = "  +5.11%  "

$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.091"
$ws.Range("E42").Value = "  +3.40%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6434"
$ws.Range("E43").Value = "  +4.67%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.190"
$ws.Range("E44").Value = "  +3.75%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.74"
$ws.Range("E46").Value = "  +5.21%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6086"
$ws.Range("E47").Value = "  +3.84%  "

$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.281"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.670"
$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.043"
$ws.Range("E50").Value = "  +6.19%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.88"
$ws.Range("E51").Value = "  +0.97%  "
